$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Password/secret.txt values between the "protected2.zip" row (row 3)
# and the "protected6.zip" row (row 7).
$b3 = $ws.Range("B3").Text
$c3 = $ws.Range("C3").Text
$b7 = $ws.Range("B7").Text
$c7 = $ws.Range("C7").Text

$ws.Range("B3").Value = $b7
$ws.Range("C3").Value = $c7
$ws.Range("B7").Value = $b3
$ws.Range("C7").Value = $c3

# Update the active cell/selection to A8
$ws.Range("A8").Select()
